# Updated symbol list on Fri Dec 30 18:54:00 UTC 2022 with GitHub Actions
# Applies the per-cell text updates to the "cryptos" sheet.
# Price-like values (column D) are assigned with a leading apostrophe so
# Excel stores them as literal text (preserving exact formatting such as
# trailing zeros) instead of re-parsing them as numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''245.10'
$ws.Range("D3").Value = '''24.96'
$ws.Range("D5").Value = '''0.05620'
$ws.Range("D7").Value = '''2.983'
$ws.Range("D8").Value = '''0.8089'
$ws.Range("D9").Value = '''0.8389'
$ws.Range("D10").Value = '''0.1338'
$ws.Range("B11").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C11").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D11").Value = '''0.03281'
$ws.Range("E11").Value = '10LiechtensteinCryptoassetsExchangeLCX'
$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D12").Value = '''0.06944'
$ws.Range("E12").Value = '11MandalaExchangeTokenMDX'
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").Value = '''0.02822'
$ws.Range("E13").Value = '12BitrueCoinBTR'
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").Value = '''0.09410'
$ws.Range("E14").Value = '13BitMartTokenBMX'
$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D15").Value = '''0.001522'
$ws.Range("E15").Value = '14BitForexTokenBF'
$ws.Range("B16").Value = 'One'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D16").Value = '''0.0005946'
$ws.Range("E16").Value = '15OneONEWorstin24h'
$ws.Range("B17").Value = 'TigerCash'
$ws.Range("C17").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D17").Value = '''0.006238'
$ws.Range("E17").Value = '16TigerCashTCH'
$ws.Range("B18").Value = 'LEO'
$ws.Range("C18").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D18").Value = '''3.500'
$ws.Range("E18").Value = '17LEOLEO'
$ws.Range("B19").Value = 'BTSEToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D19").Value = '''2.054'
$ws.Range("E19").Value = '18BTSETokenBTSE'
$ws.Range("B20").Value = 'BitpandaEcosystemToken'
$ws.Range("C20").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D20").Value = '''0.3186'
$ws.Range("E20").Value = '19BitpandaEcosystemTokenBEST'
$ws.Range("D23").Value = '''0.04688'
$ws.Range("D27").Value = '''0.00009692'
$ws.Range("D28").Value = '''0.00007231'
$ws.Range("D40").Value = '''0.03629'
$ws.Range("B41").Value = 'KickToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("D41").Value = '''0.006225'
$ws.Range("E41").Value = '40KickTokenKICK'
$ws.Range("B42").Value = 'BKEXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D42").Value = '''0.1052'
$ws.Range("E42").Value = '41BKEXTokenBKK'
$ws.Range("B43").Value = 'CEJI'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("D43").Value = '''0.002721'
$ws.Range("E43").Value = '42CEJICEJI'
$ws.Range("D44").Value = '''0.008334'
$ws.Range("D45").Value = '''0.00005278'
$ws.Range("D46").Value = '''0.00000000749'
$ws.Range("D47").Value = '''0.1899'
$ws.Range("E47").Value = '46CoinbaseStockTokenCOIN'
$ws.Range("D49").Value = '''0.00002098'
$ws.Range("D50").Value = '''0.0001998'
